# Deploying to gh-pages - update 5.4.1 sheet with City/Village/Man/Woman breakdown rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert new blank rows so the 4 existing summary rows (5,6,7,8) become
#    rows 5, 10, 15, 20, each followed by 4 new breakdown rows.
# ---------------------------------------------------------------------------
$ws.Range("A6:E9").EntireRow.Insert()      # old row6 -> row10, old row7 -> row11, old row8 -> row12
$ws.Range("A11:E14").EntireRow.Insert()    # old row7(now11) -> row15, old row8(now12) -> row16
$ws.Range("A16:E19").EntireRow.Insert()    # old row8(now16) -> row20
$ws.Range("A21:E24").EntireRow.Insert()    # room for the 4 new rows after row20

# ---------------------------------------------------------------------------
# 2. Helper to fill a row of 5 cells (A..E) with shared-string / numeric data.
# ---------------------------------------------------------------------------
function Fill-Row {
    param(
        [int]$r,
        [string]$a, [string]$b, [string]$c,
        $d, $e
    )
    $ws.Cells.Item($r, 1).Value2 = $a
    $ws.Cells.Item($r, 2).Value2 = $b
    $ws.Cells.Item($r, 3).Value2 = $c
    $ws.Cells.Item($r, 4).Value2 = $d
    $ws.Cells.Item($r, 5).Value2 = $e
}

# ---------------------------------------------------------------------------
# 3. Fill in the new breakdown rows (City / Village / Man / Woman) for each
#    of the 4 indicator categories.
# ---------------------------------------------------------------------------
# Housekeeping (row 5 summary) -> rows 6-9
Fill-Row 6  "Шаар жерлери" "Городские поселения" "City"   12.5 10.7
Fill-Row 7  "Айыл аймагы"  "Сельская местность"  "Village" 13.9 11.9
Fill-Row 8  "Эркектер"     "Мужчины"             "Man"     6.5  3.9
Fill-Row 9  "Аялдар"       "Женщины"             "Woman"   18.8 18.1

# Work on a personal/country/garden plot (row 10 summary) -> rows 11-14
Fill-Row 11 "Шаар жерлери" "Городские поселения" "City"    0.4 1
Fill-Row 12 "Айыл аймагы"  "Сельская местность"  "Village" 3.1 3.5
Fill-Row 13 "Эркектер"     "Мужчины"             "Man"     2.4 3.6
Fill-Row 14 "Аялдар"       "Женщины"             "Woman"   1   1.7

# Parenting (row 15 summary) -> rows 16-19
Fill-Row 16 "Шаар жерлери" "Городские поселения" "City"    0.9 1.7
Fill-Row 17 "Айыл аймагы"  "Сельская местность"  "Village" 0.9 2.2
Fill-Row 18 "Эркектер"     "Мужчины"             "Man"     0.6 1.2
Fill-Row 19 "Аялдар"       "Женщины"             "Woman"   1.2 2.8

# Help for relatives and friends (row 20 summary) -> rows 21-24
Fill-Row 21 "Шаар жерлери" "Городские поселения" "City"    0.4 0.3
Fill-Row 22 "Айыл аймагы"  "Сельская местность"  "Village" 0.6 0.3
Fill-Row 23 "Эркектер"     "Мужчины"             "Man"     0.6 0.3
Fill-Row 24 "Аялдар"       "Женщины"             "Woman"   0.3 0.3

# ---------------------------------------------------------------------------
# 4. Formatting: the 4 category/summary rows (5, 10, 15, 20) become bold,
#    matching the existing bold-9pt "Наименование показателей" header font
#    already used elsewhere in the sheet (row 4 / D4:E4).
# ---------------------------------------------------------------------------
$ws.Range("A5:E5").Font.Bold = $true
$ws.Range("A10:B10").Font.Bold = $true
$ws.Range("C10:E10").Font.Bold = $true
$ws.Range("A15:D15").Font.Bold = $true
$ws.Range("E15").Font.Bold = $true
$ws.Range("A20:E20").Font.Bold = $true

# Wrap text for the 2-column (Items / long names) cells of the "Work on a
# personal plot" block, which needs more vertical room (row 10 is taller).
$ws.Range("A10:B10").WrapText = $true
$ws.Range("A11:B14").WrapText = $true

# Number format 0.0 (existing custom numFmt 164) for the % columns beneath
# the Parenting block (rows 15-19, column E).
$ws.Range("E15").NumberFormat = "0.0"
$ws.Range("E16:E19").NumberFormat = "0.0"

# Row heights that carried across from the original layout.
$ws.Range("A6:E8").RowHeight = 16.5
$ws.Range("A10:E10").RowHeight = 24

# Bottom border for the final row of the table (row 24), matching the
# bottom (medium) border style used elsewhere in the sheet (row 3, which
# still carries the original untouched border style).
$ws.Range("A3:E3").Copy() | Out-Null
$ws.Range("A24:E24").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

Write-Host "Rows filled"
